$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "97.62", "0.0826") keep their exact original text representation
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply updated cell values from the latest cryptos data pull
$ws.Range("D2").Value = '42.889.16'
$ws.Range("E2").Value = '  -5.39%  '
$ws.Range("D3").Value = '2.211.28'
$ws.Range("E3").Value = '  -6.61%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '314.33'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = '97.62'
$ws.Range("E6").Value = '  -9.89%  '
$ws.Range("D7").Value = '0.580'
$ws.Range("E7").Value = '  -7.64%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -9.73%  '
$ws.Range("D10").Value = '36.44'
$ws.Range("E10").Value = '  -11.40%  '
$ws.Range("D11").Value = '54.27'
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").Value = '0.0826'
$ws.Range("E12").Value = '  -10.24%  '
$ws.Range("E13").Value = '  -8.76%  '
$ws.Range("E14").Value = '  -4.12%  '
$ws.Range("D15").Value = '0.861'
$ws.Range("E15").Value = '  -12.29%  '
$ws.Range("D16").Value = '2.547.96'
$ws.Range("E16").Value = '  -6.66%  '
$ws.Range("D17").Value = '14.06'
$ws.Range("E17").Value = '  -7.67%  '
$ws.Range("D18").Value = '2.207.49'
$ws.Range("E18").Value = '  -6.83%  '
$ws.Range("D19").Value = '42.789.94'
$ws.Range("E19").Value = '  -5.55%  '
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").Value = '0.0₃0955'
$ws.Range("E21").Value = '  -10.02%  '
$ws.Range("E22").Value = '  -12.83%  '
$ws.Range("D23").Value = '65.16'
$ws.Range("E23").Value = '  -11.01%  '
$ws.Range("D24").Value = '3.17'
$ws.Range("E24").Value = '  -9.23%  '
$ws.Range("D25").Value = '236.14'
$ws.Range("E25").Value = '  -9.38%  '
$ws.Range("E26").Value = '  -8.54%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").Value = '10.01'
$ws.Range("E28").Value = '  -10.38%  '
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  -5.55%  '
$ws.Range("D30").Value = '6.22'
$ws.Range("E30").Value = '  -15.29%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '20.41'
$ws.Range("E31").Value = '  -8.73%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.0876'
$ws.Range("E32").Value = '  -9.83%  '
$ws.Range("D33").Value = '33.70'
$ws.Range("E33").Value = '  -10.57%  '
$ws.Range("D34").Value = '154.47'
$ws.Range("E34").Value = '  -8.67%  '
$ws.Range("D35").Value = '2.78'
$ws.Range("E35").Value = '  -5.85%  '
$ws.Range("D36").Value = '3.19'
$ws.Range("E36").Value = '  +7.58%  '
$ws.Range("E37").Value = '  +14.09%  '
$ws.Range("E38").Value = '  -6.72%  '
$ws.Range("D39").Value = '4.40'
$ws.Range("E39").Value = '  -8.06%  '
$ws.Range("E40").Value = '  -12.89%  '
$ws.Range("D41").Value = '3.69'
$ws.Range("E41").Value = '  -6.07%  '
$ws.Range("E42").Value = '  -8.94%  '
$ws.Range("D43").Value = '1.858.64'
$ws.Range("E43").Value = '  +11.42%  '
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").Value = '12.28'
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("D46").Value = '88.77'
$ws.Range("E46").Value = '  -10.81%  '
$ws.Range("E47").Value = '  -11.51%  '
$ws.Range("D48").Value = '5.40'
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").Value = '75.81'
$ws.Range("E49").Value = '  -6.43%  '
$ws.Range("D50").Value = '59.89'
$ws.Range("E50").Value = '  -13.96%  '
$ws.Range("E51").Value = '  -6.48%  '
